# Refresh the cryptos price/volume snapshot table (columns B:E, rows 2-51)
# to match the latest GitHub Actions scrape. Every cell in B:E is stored as
# literal text in the sheet (inline string), including price strings that
# happen to look like plain numbers, e.g. "0.999" or "2.90".
#
# A plain `Range.Value = "0.999"` assignment lets Excel auto-coerce that
# into the Number 0.999 (and silently drops the trailing zero in a value
# like "2.90" -> 2.9), which would not match the source file's text
# formatting. So, for any column D price that parses as a plain number,
# the value is entered with a leading apostrophe -- the normal Excel
# "force text" convention -- and the cell's style is then reset to Normal
# so no stray quote-prefix formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.368.62'
$ws.Range("E2").Value = '  +6.43%  '
$ws.Range("D3").Value = '3.771.65'
$ws.Range("E3").Value = '  +21.87%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''615.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.45%  '
$ws.Range("D6").Value = '''180.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").Value = '3.768.49'
$ws.Range("E7").Value = '  +21.81%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +6.15%  '
$ws.Range("D10").Value = '''0.167'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.62%  '
$ws.Range("D11").Value = '''6.48'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  +7.41%  '
$ws.Range("D13").Value = '''40.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.50%  '
$ws.Range("D14").Value = '''0.0000258'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.48%  '
$ws.Range("D15").Value = '4.396.01'
$ws.Range("E15").Value = '  +21.74%  '
$ws.Range("D16").Value = '3.778.05'
$ws.Range("E16").Value = '  +22.06%  '
$ws.Range("D17").Value = '71.379.18'
$ws.Range("E17").Value = '  +6.50%  '
$ws.Range("D19").Value = '''7.56'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.92%  '
$ws.Range("D20").Value = '''524.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.02%  '
$ws.Range("D21").Value = '''16.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.00%  '
$ws.Range("D22").Value = '''9.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +21.51%  '
$ws.Range("D23").Value = '''0.749'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.35%  '
$ws.Range("E24").Value = '  +12.09%  '
$ws.Range("D25").Value = '''88.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.20%  '
$ws.Range("D26").Value = '''13.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.06%  '
$ws.Range("D27").Value = '''11.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.89%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '''2.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.96%  '
$ws.Range("D30").Value = '''8.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '''0.0000116'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +23.39%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.11%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''32.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.65%  '
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("D35").Value = '''0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("E36").Value = '  +12.00%  '
$ws.Range("D37").Value = '''6.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.35%  '
$ws.Range("E38").Value = '  +10.54%  '
$ws.Range("E39").Value = '  +10.60%  '
$ws.Range("D40").Value = '''0.134'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.29%  '
$ws.Range("D41").Value = '''51.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.11%  '
$ws.Range("D42").Value = '''436.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +18.43%  '
$ws.Range("D43").Value = '3.171.36'
$ws.Range("E43").Value = '  +13.18%  '
$ws.Range("D44").Value = '''8.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.64%  '
$ws.Range("D45").Value = '''44.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.04%  '
$ws.Range("D46").Value = '''2.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.58%  '
$ws.Range("E47").Value = '  +6.10%  '
$ws.Range("D48").Value = '''28.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.34%  '
$ws.Range("D49").Value = '''140.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.38%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '''2.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.89%  '
